$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entered in the same order as the original author to preserve shared-string indices

# Row 1
$ws.Range("R1").Value = "Carl"

# Row 6 (new row)
$ws.Range("R6").Value = "Normalized"
$ws.Range("T6").Value = "(Divide by 255)"

# Row 8 label first
$ws.Range("R8").Value = "Basic CNN"

# Row 7 headers
$ws.Range("T7").Value = "Acc"
$ws.Range("V7").Value = "Loss"

# Row 20/21 decision tree section
$ws.Range("R20").Value = "Decision Tree"
$ws.Range("U20").Value = "M-fold (training data only)"

# Row 9 note
$ws.Range("T9").Value = "92.48 with a random seed"

# Remaining numeric / already-shared-string cells
$ws.Range("T8").Value = 92.23
$ws.Range("V8").Value = 0.2199
$ws.Range("T20").Value = "Testing"
$ws.Range("T21").Value = 79.25
$ws.Range("U21").Value = 79.52

# Update selection to match target state
$ws.Range("U23").Select()
